$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refreshed MeteoCat extraction pass (2026-02-28 21:18-21:20 TU): new DATA_EXTRACCIO
# timestamps plus the station readings that shifted between the 20:48-20:50 and
# 21:18-21:20 snapshots. Values with a leading "'" are forced to text so Excel
# doesn't reinterpret e.g. '78%' as the number 0.78.
$ws.Cells.Item(2, 5).Value = "2026-02-28 21:18:34"
$ws.Cells.Item(2, 8).Value = "'78%"
$ws.Cells.Item(2, 15).Value = "2.8 °C"
$ws.Cells.Item(3, 5).Value = "2026-02-28 21:18:36"
$ws.Cells.Item(3, 14).Value = "-3.2 °C 20:55 TU"
$ws.Cells.Item(4, 5).Value = "2026-02-28 21:18:39"
$ws.Cells.Item(4, 15).Value = "10.9 °C"
$ws.Cells.Item(5, 5).Value = "2026-02-28 21:18:41"
$ws.Cells.Item(5, 14).Value = "-3.2 °C 20:54 TU"
$ws.Cells.Item(5, 15).Value = "-1.4 °C"
$ws.Cells.Item(6, 5).Value = "2026-02-28 21:18:43"
$ws.Cells.Item(6, 14).Value = "9.1 °C 20:46 TU"
$ws.Cells.Item(6, 15).Value = "11.9 °C"
$ws.Cells.Item(7, 5).Value = "2026-02-28 21:18:45"
$ws.Cells.Item(8, 5).Value = "2026-02-28 21:18:48"
$ws.Cells.Item(8, 10).Value = "1024.8 hPa"
$ws.Cells.Item(9, 5).Value = "2026-02-28 21:18:51"
$ws.Cells.Item(10, 5).Value = "2026-02-28 21:18:54"
$ws.Cells.Item(10, 15).Value = "10.7 °C"
$ws.Cells.Item(11, 5).Value = "2026-02-28 21:18:56"
$ws.Cells.Item(12, 5).Value = "2026-02-28 21:18:59"
$ws.Cells.Item(13, 5).Value = "2026-02-28 21:19:01"
$ws.Cells.Item(13, 11).Value = "11.9 MJ/m2"
$ws.Cells.Item(14, 5).Value = "2026-02-28 21:19:04"
$ws.Cells.Item(15, 5).Value = "2026-02-28 21:19:06"
$ws.Cells.Item(15, 8).Value = "'81%"
$ws.Cells.Item(15, 15).Value = "10.9 °C"
$ws.Cells.Item(16, 5).Value = "2026-02-28 21:19:09"
$ws.Cells.Item(16, 14).Value = "-3.1 °C 20:32 TU"
$ws.Cells.Item(17, 5).Value = "2026-02-28 21:19:12"
$ws.Cells.Item(17, 14).Value = "1.4 °C 20:30 TU"
$ws.Cells.Item(18, 5).Value = "2026-02-28 21:19:14"
$ws.Cells.Item(18, 10).Value = "1025.0 hPa"
$ws.Cells.Item(18, 15).Value = "11.5 °C"
$ws.Cells.Item(19, 5).Value = "2026-02-28 21:19:17"
$ws.Cells.Item(19, 15).Value = "7.8 °C"
$ws.Cells.Item(20, 5).Value = "2026-02-28 21:19:20"
$ws.Cells.Item(20, 8).Value = "'65%"
$ws.Cells.Item(20, 14).Value = "-2.2 °C 20:40 TU"
$ws.Cells.Item(20, 15).Value = "-0.6 °C"
$ws.Cells.Item(21, 5).Value = "2026-02-28 21:19:22"
$ws.Cells.Item(21, 10).Value = "1024.1 hPa"
$ws.Cells.Item(22, 5).Value = "2026-02-28 21:19:25"
$ws.Cells.Item(22, 8).Value = "'70%"
$ws.Cells.Item(22, 14).Value = "-3.1 °C 20:39 TU"
$ws.Cells.Item(22, 15).Value = "-1.7 °C"
$ws.Cells.Item(23, 5).Value = "2026-02-28 21:19:28"
$ws.Cells.Item(23, 8).Value = "'72%"
$ws.Cells.Item(23, 9).Value = "1.6 mm"
$ws.Cells.Item(24, 5).Value = "2026-02-28 21:19:30"
$ws.Cells.Item(24, 15).Value = "8.5 °C"
$ws.Cells.Item(25, 5).Value = "2026-02-28 21:19:32"
$ws.Cells.Item(25, 14).Value = "-1.1 °C 20:46 TU"
$ws.Cells.Item(25, 15).Value = "1.1 °C"
$ws.Cells.Item(26, 5).Value = "2026-02-28 21:19:35"
$ws.Cells.Item(27, 5).Value = "2026-02-28 21:19:37"
$ws.Cells.Item(27, 8).Value = "'57%"
$ws.Cells.Item(27, 14).Value = "-0.6 °C 20:59 TU"
$ws.Cells.Item(28, 5).Value = "2026-02-28 21:19:40"
$ws.Cells.Item(28, 10).Value = "1024.8 hPa"
$ws.Cells.Item(29, 5).Value = "2026-02-28 21:19:43"
$ws.Cells.Item(30, 5).Value = "2026-02-28 21:19:46"
$ws.Cells.Item(31, 5).Value = "2026-02-28 21:19:48"
$ws.Cells.Item(31, 8).Value = "'80%"
$ws.Cells.Item(31, 12).Value = "68.4 km/h - 347º 20:40 TU"
$ws.Cells.Item(32, 5).Value = "2026-02-28 21:19:51"
$ws.Cells.Item(33, 5).Value = "2026-02-28 21:19:54"
$ws.Cells.Item(34, 5).Value = "2026-02-28 21:19:57"
$ws.Cells.Item(34, 8).Value = "'69%"
$ws.Cells.Item(34, 9).Value = "1.1 mm"
$ws.Cells.Item(34, 15).Value = "1.2 °C"
$ws.Cells.Item(35, 5).Value = "2026-02-28 21:20:00"
$ws.Cells.Item(36, 5).Value = "2026-02-28 21:20:02"
$ws.Cells.Item(36, 15).Value = "12.6 °C"
$ws.Cells.Item(37, 5).Value = "2026-02-28 21:20:05"
$ws.Cells.Item(37, 10).Value = "1025.9 hPa"
$ws.Cells.Item(38, 5).Value = "2026-02-28 21:20:08"
$ws.Cells.Item(38, 8).Value = "'81%"
$ws.Cells.Item(38, 15).Value = "11.6 °C"
$ws.Cells.Item(39, 5).Value = "2026-02-28 21:20:10"
$ws.Cells.Item(40, 5).Value = "2026-02-28 21:20:13"
$ws.Cells.Item(41, 5).Value = "2026-02-28 21:20:15"
$ws.Cells.Item(41, 15).Value = "13.2 °C"
$ws.Cells.Item(42, 5).Value = "2026-02-28 21:20:18"
$ws.Cells.Item(42, 8).Value = "'88%"
$ws.Cells.Item(42, 15).Value = "11.0 °C"
$ws.Cells.Item(43, 5).Value = "2026-02-28 21:20:20"
$ws.Cells.Item(43, 8).Value = "'81%"
$ws.Cells.Item(43, 15).Value = "7.8 °C"
$ws.Cells.Item(44, 5).Value = "2026-02-28 21:20:23"
$ws.Cells.Item(44, 14).Value = "-3.0 °C 20:38 TU"
$ws.Cells.Item(44, 15).Value = "-1.2 °C"
$ws.Cells.Item(45, 5).Value = "2026-02-28 21:20:26"
$ws.Cells.Item(45, 10).Value = "1025.6 hPa"
$ws.Cells.Item(45, 14).Value = "4.2 °C 20:58 TU"
$ws.Cells.Item(46, 5).Value = "2026-02-28 21:20:28"
$ws.Cells.Item(46, 8).Value = "'79%"
$ws.Cells.Item(46, 15).Value = "11.4 °C"
